# Update "想去人数" (interested-count) values in column F across the
# four worksheets of the 广州-漫展信息 workbook, per the regenerated
# gh-pages data snapshot (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(3, 6).Value  = 1402
$ws1.Cells.Item(4, 6).Value  = 26233
$ws1.Cells.Item(7, 6).Value  = 589
$ws1.Cells.Item(8, 6).Value  = 171
$ws1.Cells.Item(9, 6).Value  = 526
$ws1.Cells.Item(11, 6).Value = 348
$ws1.Cells.Item(12, 6).Value = 214
$ws1.Cells.Item(15, 6).Value = 289
$ws1.Cells.Item(17, 6).Value = 363
$ws1.Cells.Item(18, 6).Value = 52
$ws1.Cells.Item(19, 6).Value = 1495
$ws1.Cells.Item(20, 6).Value = 174
$ws1.Cells.Item(21, 6).Value = 22

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(3, 6).Value  = 224
$ws2.Cells.Item(6, 6).Value  = 166
$ws2.Cells.Item(8, 6).Value  = 108
$ws2.Cells.Item(9, 6).Value  = 108
$ws2.Cells.Item(10, 6).Value = 431
$ws2.Cells.Item(15, 6).Value = 49
$ws2.Cells.Item(19, 6).Value = 19

# Sheet "本地生活"
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Cells.Item(2, 6).Value = 4964
$ws3.Cells.Item(3, 6).Value = 202

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(3, 6).Value  = 1402
$ws4.Cells.Item(4, 6).Value  = 4964
$ws4.Cells.Item(5, 6).Value  = 202
$ws4.Cells.Item(6, 6).Value  = 26233
$ws4.Cells.Item(10, 6).Value = 224
$ws4.Cells.Item(11, 6).Value = 589
$ws4.Cells.Item(14, 6).Value = 171
$ws4.Cells.Item(15, 6).Value = 166
$ws4.Cells.Item(16, 6).Value = 166
$ws4.Cells.Item(18, 6).Value = 108
$ws4.Cells.Item(19, 6).Value = 108
$ws4.Cells.Item(20, 6).Value = 431
$ws4.Cells.Item(21, 6).Value = 526
$ws4.Cells.Item(24, 6).Value = 348
$ws4.Cells.Item(25, 6).Value = 214
$ws4.Cells.Item(29, 6).Value = 289
$ws4.Cells.Item(33, 6).Value = 363
$ws4.Cells.Item(34, 6).Value = 52
$ws4.Cells.Item(35, 6).Value = 49
$ws4.Cells.Item(36, 6).Value = 1495
$ws4.Cells.Item(37, 6).Value = 174
$ws4.Cells.Item(39, 6).Value = 22
$ws4.Cells.Item(45, 6).Value = 19
